# Update calibration figures across the "PnL" calibration sheets for the
# finalized virtual-environment definitions.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.7105803991164821
$ws.Range("B3").Value = 0.8072336039972227
$ws.Range("B4").Value = 64654.8047545133
$ws.Range("B5").Value = "[1.0, 0.33270968067201345, 0.3522734129625214, 0.32317918370065557, 0.30751957730981044, 0.5162282386426504, 0.2100249284164006, 0.26145819483095534, 0.22241639251953693, 0.21618163908983487, 0.2097338444319448, 0.22983995957359743, 0.24028644478531952, 0.2352689579679485, 0.23853566729011302, 0.229227694674077, 0.24456981283717785, 0.23185946550435874, 0.21572016316242312, 0.23334296872368138]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.8918518518518519
$ws.Range("B4").Value = -34.66827684399614
$ws.Range("B5").Value = 0.8103193911293781
$ws.Range("B6").Value = 88069.38718990301
$ws.Range("B7").Value = 30.25273433222138
$ws.Range("B8").Value = 0.8060331275414746
$ws.Range("B9").Value = 23924.86992146928

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.003540348544720932
$ws.Range("B3").Value = 0.224553365954934
$ws.Range("B4").Value = 0.1664891910778499
$ws.Range("B5").Value = 0.8335108089221501
$ws.Range("B6").Value = "[1.0, -0.06389447414691769, -0.021975633464639367, -0.05575591567177519, -0.09704706102532343, 0.24467287959117287, -0.028392294169012077, 0.04199393959522337, 0.02389219238773284, -0.044520016612905444, -0.053419806800431444, -0.04249877712043097, -0.0024519557944264435, 0.030206687791300216, 0.02740347105342912, 0.017151476334452184, -0.016246990130471456, -0.022908853398877122, -0.02460583454380923, 0.006163781898280912]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.08574797181727643
$ws.Range("B3").Value = 0.2523070119632156
$ws.Range("B4").Value = 0.1095130545210253
$ws.Range("B5").Value = 0.8282956973118958
$ws.Range("B6").Value = "[0.9999999999999999, -0.06165408236972, -0.02102261560184984, -0.05412837581694195, -0.09478804476315274, 0.24764344111519113, -0.030799598786885055, 0.038786089231267334, 0.022953628188991482, -0.04866808902365025, -0.05570819772229042, -0.044411233414498895, -0.0033943978177947247, 0.029564546165195608, 0.0258951925090044, 0.017217190832167495, -0.016745238004780196, -0.021804353575299584, -0.02323154498720819, 0.008421973393838664]"
$ws.Range("B7").Value = 0.1244339196633769

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.00503610898272411
$ws.Range("B3").Value = 0.1775498965660597
$ws.Range("B4").Value = 0.1641262696287309
$ws.Range("B5").Value = 0.8476452660947248
$ws.Range("B6").Value = "[1.0, -0.03593281601307041, -0.010639861439722571, -0.03951461704850044, -0.05710298676966572, 0.17424670563482386, -0.027647724843411707, 0.02627153565433829, 0.006657593043323032, -0.044100236783981464, -0.03754688623655121, -0.03250097729821977, -0.016952587219524567, 0.03712884077282447, 0.01965000872545261, 0.00751816919851212, -0.008385327153204718, -0.020185537453967596, -0.009379488419193392, 0.00645869935094012]"
$ws.Range("B7").Value = -0.0235005591541813
$ws.Range("B9").Value = 0.7967641922719754
